# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values computed for column G (rows 2-21)
$sVals = @{
    2  = 1
    3  = 5
    4  = 4
    5  = 7
    6  = 6
    7  = 4
    8  = 5
    9  = 3
    10 = 7
    11 = 10
    12 = 4
    13 = 6
    14 = 3
    15 = 6
    16 = 7
    17 = 2
    18 = 5
    19 = 2
    20 = 1
    21 = 2
}

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $sVals[$row]
}
